# Update: pi 15. 01. 2021
# Revises historical AgTests (H) / AgPosit (I) cumulative columns and
# appends the new daily row (316) for 2021-01-15 (date serial 44210).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Revised AgTests (H) / AgPosit (I) values for existing rows ---
$updates = @(
    @{ Row = 286; H = 54344; I = 4260 },
    @{ Row = 287; H = 57881; I = 3937 },
    @{ Row = 288; H = 56687; I = 3993 },
    @{ Row = 289; H = 65263; I = 3764 },
    @{ Row = 292; H = 82425; I = 7311 },
    @{ Row = 293; H = 83197; I = 5869 },
    @{ Row = 294; H = 92254; I = 5110 },
    @{ Row = 298; H = 3048;  I = 281 },
    @{ Row = 299; H = 65636; I = 6884 },
    @{ Row = 300; H = 71201; I = 6963 },
    @{ Row = 301; H = 70185; I = 5566 },
    @{ Row = 302; H = 73183; I = 5336 },
    @{ Row = 307; H = 73167; I = 6329 },
    @{ Row = 309; H = 57732; I = 4014 },
    @{ Row = 310; H = 90050; I = 5389 },
    @{ Row = 313; H = 72406; I = 3515 },
    @{ Row = 314; H = 64130; I = 3299 },
    @{ Row = 315; H = 64594; I = 3034 }
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 8).Value = $u.H
    $ws.Cells.Item($u.Row, 9).Value = $u.I
}

# --- Append new row 316 (2021-01-15) ---
$newRow = 316
$ws.Cells.Item($newRow, 1).Value = 44210
$ws.Cells.Item($newRow, 1).NumberFormat = "yyyy-mm-dd"
$ws.Cells.Item($newRow, 2).Value = 220707
$ws.Cells.Item($newRow, 3).Value = 163323
$ws.Cells.Item($newRow, 4).Value = 54022
$ws.Cells.Item($newRow, 5).Value = 11875
$ws.Cells.Item($newRow, 6).Value = 2729
$ws.Cells.Item($newRow, 7).Value = 3362
$ws.Cells.Item($newRow, 8).Value = 46169
$ws.Cells.Item($newRow, 9).Value = 2165
